$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was "Blåsippa", becomes "Tretåig hackspett" (values rotated from old row 6)
$ws.Range("A3").Value = 112042452
$ws.Range("B3").Value = 56430
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("P3").Value = "Stor Mpmerg, Kilen-Stor, Moberg, Leksand, Dlr"
$ws.Range("Q3").Value = 511614
$ws.Range("R3").Value = 6733640
$ws.Range("S3").Value = 25
$ws.Range("Z3").Value = ""
$ws.Range("AB3").Value = ""
$ws.Range("AC3").Value = ""
$ws.Range("AW3").Value = "Åke Sköld"
$ws.Range("AX3").Value = "Åke Sköld"

# Row 4: minor update
$ws.Range("B4").Value = 90480

# Row 5: was "Revlummer", becomes "Blåsippa" (values rotated from old row 3)
$ws.Range("A5").Value = 112042940
$ws.Range("B5").Value = 98980
$ws.Range("E5").Value = 222498
$ws.Range("F5").Value = "Blåsippa"
$ws.Range("G5").Value = "Hepatica nobilis"
$ws.Range("H5").Value = "Schreb."
$ws.Range("Q5").Value = 511611
$ws.Range("R5").Value = 6733626
$ws.Range("Z5").Value = "10:33"
$ws.Range("AB5").Value = "10:33"
$ws.Range("AC5").Value = "Fullt med blåsippsblad på denna sidan bäcken"

# Row 6: was "Tretåig hackspett", becomes "Revlummer" (values rotated from old row 5)
$ws.Range("A6").Value = 112043158
$ws.Range("B6").Value = 95701
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221945
$ws.Range("F6").Value = "Revlummer"
$ws.Range("G6").Value = "Lycopodium annotinum"
$ws.Range("H6").Value = "L."
$ws.Range("P6").Value = "Stor-Moberg (Stor-Moberg), Dlr"
$ws.Range("Q6").Value = 511628
$ws.Range("R6").Value = 6733623
$ws.Range("S6").Value = 1
$ws.Range("Z6").Value = "10:51"
$ws.Range("AB6").Value = "10:51"
$ws.Range("AC6").Value = "Finns fläckvis i området"
$ws.Range("AW6").Value = "Evalena Sköld"
$ws.Range("AX6").Value = "Evalena Sköld, Åke Sköld"
